$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 24 for Exp 27 experiment parameters (enter the label first)
$ws.Range("A24").Value = "Exp 27"
$ws.Range("B24").Value = 0.7
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = "Local"
$ws.Range("E24").Value = -1

# Fill in the missing "Result Image Name" for Exp 26 (row 23)
$ws.Range("F23").Value = "Exp 26.png"

# Result Image Name for the new Exp 27 row
$ws.Range("F24").Value = "Exp 27.png"

# Match cell styling (centered alignment) used by the rest of the data rows
$ws.Range("A24:E24").HorizontalAlignment = -4108

# Update view state: scroll position and active selection after edits
$ws.Range("F25").Select()
$excel.ActiveWindow.ScrollRow = 6
